$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.376.72'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '1.874.17'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7115'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3117'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07787'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08466'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('D12').Value = '1.873.29'
$ws.Range('E12').Value = '  -4.66%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.240'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7128'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.25'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('D16').Value = '29.380.76'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008246'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.80%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.041'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '240.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.26'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').Value = '2.125.03'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.798'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1607'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '163.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.077'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.510'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.420'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.279'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.64%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05308'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.32%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.936'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.178'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.45%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7486'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.93%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.696'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01872'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.721'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('D40').Value = '1.204.78'
$ws.Range('E40').Value = '  +1.90%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.447'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '73.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8877'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '107.76'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.31%  '
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '2.021.42'
$ws.Range('E46').Value = '  -2.15%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.820'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5209'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('E49').Value = '  +7.78%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.397'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('E51').Value = '  +0.86%  '
